$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update F-column counts
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 9026
$wsExhibit.Range("F10").Value = 825
$wsExhibit.Range("F13").Value = 1028
$wsExhibit.Range("F14").Value = 127
$wsExhibit.Range("F18").Value = 313
$wsExhibit.Range("F20").Value = 239
$wsExhibit.Range("F21").Value = 1150

# Sheet "全部类型" (All types) - update corresponding F-column counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 9026
$wsAll.Range("F12").Value = 825
$wsAll.Range("F15").Value = 1028
$wsAll.Range("F16").Value = 127
$wsAll.Range("F20").Value = 313
$wsAll.Range("F22").Value = 239
$wsAll.Range("F23").Value = 1150
